# Applies the 2026-02-01 20:48:26 scrape refresh to the "horarios-141" workbook.
# Generated from a precise cell-level diff between the before/after OOXML.
$wb = $excel.ActiveWorkbook

# ---- Worksheet 1: LP1912 ----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 1).Value = 'Última actualización: 20:48:26'
$ws.Cells.Item(3, 1).Value = 'Total filas: 320'
$ws.Cells.Item(66, 1).Value = '08:49:06'
$ws.Cells.Item(66, 3).Value = '14_ABASTO'
$ws.Cells.Item(66, 4).Value = 29
$ws.Cells.Item(67, 1).Value = '08:57:42'
$ws.Cells.Item(67, 3).Value = '15X38_ABASTO'
$ws.Cells.Item(67, 4).Value = 21
$ws.Cells.Item(90, 3).Value = '15_ABASTO'
$ws.Cells.Item(91, 3).Value = '14_ABASTO'
$ws.Cells.Item(103, 1).Value = '11:01:19'
$ws.Cells.Item(103, 3).Value = '15X38_ABASTO'
$ws.Cells.Item(103, 4).Value = 14
$ws.Cells.Item(104, 1).Value = '10:32:07'
$ws.Cells.Item(104, 3).Value = '14_ABASTO'
$ws.Cells.Item(104, 4).Value = 43
$ws.Cells.Item(117, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(118, 3).Value = '17_ROMERO'
$ws.Cells.Item(125, 3).Value = '10_OLMOS'
$ws.Cells.Item(126, 3).Value = '215C_EL PATO'
$ws.Cells.Item(133, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(134, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(136, 1).Value = '12:43:13'
$ws.Cells.Item(136, 3).Value = '14_ABASTO'
$ws.Cells.Item(136, 4).Value = 5
$ws.Cells.Item(137, 3).Value = '15X38_ABASTO'
$ws.Cells.Item(138, 1).Value = '12:18:38'
$ws.Cells.Item(138, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(138, 4).Value = 30
$ws.Cells.Item(140, 1).Value = '12:18:38'
$ws.Cells.Item(140, 3).Value = '215C_EL PATO'
$ws.Cells.Item(140, 4).Value = 45
$ws.Cells.Item(141, 1).Value = '12:43:13'
$ws.Cells.Item(141, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(141, 4).Value = 20
$ws.Cells.Item(151, 1).Value = '13:28:27'
$ws.Cells.Item(151, 3).Value = '215A_EL PATO'
$ws.Cells.Item(151, 4).Value = 5
$ws.Cells.Item(152, 1).Value = '12:43:13'
$ws.Cells.Item(152, 3).Value = '14_ABASTO'
$ws.Cells.Item(152, 4).Value = 50
$ws.Cells.Item(187, 1).Value = '15:34:15'
$ws.Cells.Item(187, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(187, 4).Value = 7
$ws.Cells.Item(188, 1).Value = '14:58:38'
$ws.Cells.Item(188, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(188, 4).Value = 43
$ws.Cells.Item(199, 1).Value = '14:58:38'
$ws.Cells.Item(199, 3).Value = '14_ABASTO'
$ws.Cells.Item(199, 4).Value = 67
$ws.Cells.Item(200, 1).Value = '16:02:30'
$ws.Cells.Item(200, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(200, 4).Value = 3
$ws.Cells.Item(230, 1).Value = '17:22:11'
$ws.Cells.Item(230, 3).Value = '17_ROMERO'
$ws.Cells.Item(230, 4).Value = 12
$ws.Cells.Item(231, 1).Value = '16:34:05'
$ws.Cells.Item(231, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(231, 4).Value = 60
$ws.Cells.Item(280, 1).Value = '18:01:05'
$ws.Cells.Item(280, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(280, 4).Value = 82
$ws.Cells.Item(281, 1).Value = '19:14:15'
$ws.Cells.Item(281, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(281, 4).Value = 9
$ws.Cells.Item(303, 1).Value = '20:48:26'
$ws.Cells.Item(303, 4).Value = 0
$ws.Cells.Item(304, 1).Value = '20:48:26'
$ws.Cells.Item(304, 2).Value = '20:49'
$ws.Cells.Item(304, 4).Value = 1
$ws.Cells.Item(305, 1).Value = '19:14:15'
$ws.Cells.Item(305, 2).Value = '20:50'
$ws.Cells.Item(305, 4).Value = 96
$ws.Cells.Item(306, 1).Value = '19:45:00'
$ws.Cells.Item(306, 2).Value = '20:52'
$ws.Cells.Item(306, 3).Value = '17_ROMERO'
$ws.Cells.Item(306, 4).Value = 67
$ws.Cells.Item(307, 1).Value = '20:48:26'
$ws.Cells.Item(307, 2).Value = '20:55'
$ws.Cells.Item(307, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(307, 4).Value = 7
$ws.Cells.Item(308, 1).Value = '20:48:26'
$ws.Cells.Item(308, 2).Value = '20:56'
$ws.Cells.Item(308, 4).Value = 8
$ws.Cells.Item(309, 2).Value = '20:57'
$ws.Cells.Item(309, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(309, 4).Value = 29
$ws.Cells.Item(310, 1).Value = '20:48:26'
$ws.Cells.Item(310, 2).Value = '21:07'
$ws.Cells.Item(310, 3).Value = '10_OLMOS'
$ws.Cells.Item(310, 4).Value = 19
$ws.Cells.Item(311, 1).Value = '20:48:26'
$ws.Cells.Item(311, 2).Value = '21:10'
$ws.Cells.Item(311, 3).Value = '15_ABASTO'
$ws.Cells.Item(311, 4).Value = 22
$ws.Cells.Item(312, 1).Value = '20:48:26'
$ws.Cells.Item(312, 2).Value = '21:28'
$ws.Cells.Item(312, 4).Value = 40
$ws.Cells.Item(313, 2).Value = '21:29'
$ws.Cells.Item(313, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(313, 4).Value = 89
$ws.Cells.Item(314, 1).Value = '20:00:07'
$ws.Cells.Item(314, 2).Value = '21:33'
$ws.Cells.Item(314, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(314, 4).Value = 93
$ws.Cells.Item(315, 1).Value = '20:48:26'
$ws.Cells.Item(315, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(315, 4).Value = 46
$ws.Cells.Item(316, 1).Value = '20:48:26'
$ws.Cells.Item(316, 2).Value = '21:34'
$ws.Cells.Item(316, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(316, 4).Value = 46
$ws.Cells.Item(317, 2).Value = '21:45'
$ws.Cells.Item(317, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(317, 4).Value = 77
$ws.Cells.Item(318, 1).Value = '20:48:26'
$ws.Cells.Item(318, 2).Value = '21:46'
$ws.Cells.Item(318, 3).Value = '14X44_ABASTO'
$ws.Cells.Item(318, 4).Value = 58
$ws.Cells.Item(319, 1).Value = '20:00:07'
$ws.Cells.Item(319, 2).Value = '21:48'
$ws.Cells.Item(319, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(319, 4).Value = 108
$ws.Cells.Item(320, 1).Value = '20:48:26'
$ws.Cells.Item(320, 2).Value = '21:55'
$ws.Cells.Item(320, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(320, 4).Value = 67
$ws.Cells.Item(321, 1).Value = '20:48:26'
$ws.Cells.Item(321, 2).Value = '22:04'
$ws.Cells.Item(321, 3).Value = '15_ABASTO'
$ws.Cells.Item(321, 4).Value = 76
$ws.Cells.Item(321, 5).Value = 'LP1912'
$ws.Cells.Item(322, 1).Value = '20:48:26'
$ws.Cells.Item(322, 2).Value = '22:11'
$ws.Cells.Item(322, 3).Value = '14_ABASTO'
$ws.Cells.Item(322, 4).Value = 83
$ws.Cells.Item(322, 5).Value = 'LP1912'
$ws.Cells.Item(323, 1).Value = '20:48:26'
$ws.Cells.Item(323, 2).Value = '22:33'
$ws.Cells.Item(323, 3).Value = '215C_EL PATO'
$ws.Cells.Item(323, 4).Value = 105
$ws.Cells.Item(323, 5).Value = 'LP1912'
$ws.Cells.Item(324, 1).Value = '20:48:26'
$ws.Cells.Item(324, 2).Value = '22:34'
$ws.Cells.Item(324, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(324, 4).Value = 106
$ws.Cells.Item(324, 5).Value = 'LP1912'
$ws.Cells.Item(325, 1).Value = '20:48:26'
$ws.Cells.Item(325, 2).Value = '22:43'
$ws.Cells.Item(325, 3).Value = '11X44_ETCHEVERRY'
$ws.Cells.Item(325, 4).Value = 115
$ws.Cells.Item(325, 5).Value = 'LP1912'

# ---- Worksheet 2: LP1912-215 ----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 1).Value = 'Última actualización: 20:48:26'
$ws.Cells.Item(3, 1).Value = 'Total filas: 49'
$ws.Cells.Item(53, 1).Value = '20:48:26'
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(54, 1).Value = '20:48:26'
$ws.Cells.Item(54, 2).Value = '22:33'
$ws.Cells.Item(54, 3).Value = '215C_EL PATO'
$ws.Cells.Item(54, 4).Value = 105
$ws.Cells.Item(54, 5).Value = 'LP1912'

# ---- Worksheet 3: 6203-6173 ----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 1).Value = 'Última actualización: 20:48:26'
$ws.Cells.Item(3, 1).Value = 'Total filas: 45'
$ws.Cells.Item(19, 1).Value = '09:42:42'
$ws.Cells.Item(19, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(19, 4).Value = 48
$ws.Cells.Item(20, 1).Value = '08:49:06'
$ws.Cells.Item(20, 3).Value = '215B_LP-P MOR-1 Y 57'
$ws.Cells.Item(20, 4).Value = 101
$ws.Cells.Item(48, 1).Value = '20:48:26'
$ws.Cells.Item(48, 2).Value = '20:54'
$ws.Cells.Item(48, 4).Value = 6
$ws.Cells.Item(49, 1).Value = '20:48:26'
$ws.Cells.Item(49, 2).Value = '21:30'
$ws.Cells.Item(49, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(49, 4).Value = 42
$ws.Cells.Item(49, 5).Value = 'L6203'
$ws.Cells.Item(50, 1).Value = '20:48:26'
$ws.Cells.Item(50, 2).Value = '22:20'
$ws.Cells.Item(50, 3).Value = '215B_LP-P MOR-40 Y 115'
$ws.Cells.Item(50, 4).Value = 92
$ws.Cells.Item(50, 5).Value = 'L6173'
